$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 1.085893666666667
$ws.Cells.Item(2, 8).Value = 3.257681
$ws.Cells.Item(2, 9).Value = 0.02840400986010362
$ws.Cells.Item(2, 10).Value = 0.02840400986010361
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 82.83048866666665
$ws.Cells.Item(2, 14).Value = 248.491466
$ws.Cells.Item(2, 15).Value = 0.3167437020391103
$ws.Cells.Item(2, 16).Value = 0.3167437020391103
$ws.Cells.Item(2, 17).Value = 89.94510305003844
$ws.Cells.Item(2, 18).Value = 809.505927450346
$ws.Cells.Item(2, 19).Value = 0.008996791235844612
$ws.Cells.Item(2, 20).Value = 0.008996791235844612

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 1.085893666666667
$ws.Cells.Item(3, 8).Value = 3.257681
$ws.Cells.Item(3, 9).Value = 0.02840400986010362
$ws.Cells.Item(3, 10).Value = 0.02840400986010361
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 152.851481
$ws.Cells.Item(3, 14).Value = 458.554443
$ws.Cells.Item(3, 15).Value = 0.5845039034954311
$ws.Cells.Item(3, 16).Value = 0.5845039034954312
$ws.Cells.Item(3, 17).Value = 165.9804551585204
$ws.Cells.Item(3, 18).Value = 1493.824096426683
$ws.Cells.Item(3, 19).Value = 0.01660225463815328
$ws.Cells.Item(3, 20).Value = 0.01660225463815328

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 1.085893666666667
$ws.Cells.Item(4, 8).Value = 3.257681
$ws.Cells.Item(4, 9).Value = 0.02840400986010362
$ws.Cells.Item(4, 10).Value = 0.02840400986010361
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 25.824378
$ws.Cells.Item(4, 14).Value = 77.47313399999999
$ws.Cells.Item(4, 15).Value = 0.09875239446545848
$ws.Cells.Item(4, 16).Value = 0.0987523944654585
$ws.Cells.Item(4, 17).Value = 28.042528515806
$ws.Cells.Item(4, 18).Value = 252.382756642254
$ws.Cells.Item(4, 19).Value = 0.002804963986105725
$ws.Cells.Item(4, 20).Value = 0.002804963986105725

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 33.91722166666667
$ws.Cells.Item(5, 8).Value = 101.751665
$ws.Cells.Item(5, 9).Value = 0.887181800778517
$ws.Cells.Item(5, 10).Value = 0.8871818007785169
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 82.83048866666665
$ws.Cells.Item(5, 14).Value = 248.491466
$ws.Cells.Item(5, 15).Value = 0.3167437020391103
$ws.Cells.Item(5, 16).Value = 0.3167437020391103
$ws.Cells.Item(5, 17).Value = 2809.380044865654
$ws.Cells.Item(5, 18).Value = 25284.42040379089
$ws.Cells.Item(5, 19).Value = 0.2810092479603119
$ws.Cells.Item(5, 20).Value = 0.2810092479603119

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 33.91722166666667
$ws.Cells.Item(6, 8).Value = 101.751665
$ws.Cells.Item(6, 9).Value = 0.887181800778517
$ws.Cells.Item(6, 10).Value = 0.8871818007785169
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 152.851481
$ws.Cells.Item(6, 14).Value = 458.554443
$ws.Cells.Item(6, 15).Value = 0.5845039034954311
$ws.Cells.Item(6, 16).Value = 0.5845039034954312
$ws.Cells.Item(6, 17).Value = 5184.297563155289
$ws.Cells.Item(6, 18).Value = 46658.6780683976
$ws.Cells.Item(6, 19).Value = 0.5185612256651492
$ws.Cells.Item(6, 20).Value = 0.5185612256651492

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 33.91722166666667
$ws.Cells.Item(7, 8).Value = 101.751665
$ws.Cells.Item(7, 9).Value = 0.887181800778517
$ws.Cells.Item(7, 10).Value = 0.8871818007785169
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 25.824378
$ws.Cells.Item(7, 14).Value = 77.47313399999999
$ws.Cells.Item(7, 15).Value = 0.09875239446545848
$ws.Cells.Item(7, 16).Value = 0.0987523944654585
$ws.Cells.Item(7, 17).Value = 875.89115302979
$ws.Cells.Item(7, 18).Value = 7883.020377268109
$ws.Cells.Item(7, 19).Value = 0.08761132715305592
$ws.Cells.Item(7, 20).Value = 0.08761132715305592

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 3.227179333333333
$ws.Cells.Item(8, 8).Value = 9.681538
$ws.Cells.Item(8, 9).Value = 0.08441418936137941
$ws.Cells.Item(8, 10).Value = 0.0844141893613794
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 82.83048866666665
$ws.Cells.Item(8, 14).Value = 248.491466
$ws.Cells.Item(8, 15).Value = 0.3167437020391103
$ws.Cells.Item(8, 16).Value = 0.3167437020391103
$ws.Cells.Item(8, 17).Value = 267.3088411949675
$ws.Cells.Item(8, 18).Value = 2405.779570754707
$ws.Cells.Item(8, 19).Value = 0.0267376628429538
$ws.Cells.Item(8, 20).Value = 0.0267376628429538

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 3.227179333333333
$ws.Cells.Item(9, 8).Value = 9.681538
$ws.Cells.Item(9, 9).Value = 0.08441418936137941
$ws.Cells.Item(9, 10).Value = 0.0844141893613794
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 152.851481
$ws.Cells.Item(9, 14).Value = 458.554443
$ws.Cells.Item(9, 15).Value = 0.5845039034954311
$ws.Cells.Item(9, 16).Value = 0.5845039034954312
$ws.Cells.Item(9, 17).Value = 493.2791405525927
$ws.Cells.Item(9, 18).Value = 4439.512264973334
$ws.Cells.Item(9, 19).Value = 0.04934042319212876
$ws.Cells.Item(9, 20).Value = 0.04934042319212877

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 3.227179333333333
$ws.Cells.Item(10, 8).Value = 9.681538
$ws.Cells.Item(10, 9).Value = 0.08441418936137941
$ws.Cells.Item(10, 10).Value = 0.0844141893613794
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 25.824378
$ws.Cells.Item(10, 14).Value = 77.47313399999999
$ws.Cells.Item(10, 15).Value = 0.09875239446545848
$ws.Cells.Item(10, 16).Value = 0.0987523944654585
$ws.Cells.Item(10, 17).Value = 83.33989897778798
$ws.Cells.Item(10, 18).Value = 750.0590908000919
$ws.Cells.Item(10, 19).Value = 0.008336103326296848
$ws.Cells.Item(10, 20).Value = 0.008336103326296848
